$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1:B3")

# Re-apply the italic purple font used for the exchange-rate code/rate block.
$rng.Font.Name = "Calibri"
$rng.Font.Size = 11
$rng.Font.Italic = $true
$rng.Font.ColorIndex = 13

# A1 holds the account code; keep it text (not auto-converted to a number)
# and update it to the new code "4004".
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "4004"
